$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last refreshed" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 13:52"

# --- Switzerland (Suiza, row 20): refresh active/recovered counts ---
$ws.Cells.Item(20, 4).Value = 23100
$ws.Cells.Item(20, 5).Value = 4770

# --- Rows 57-61: Kuwait moves up (alphabetically/by rank) ahead of
#     Argelia/Moldavia/Luxemburgo, each of those shifting down one row,
#     and Kazajistan (row 61) gets refreshed figures ---
$ws.Cells.Item(57, 1).Value = "Kuwait"
$ws.Cells.Item(57, 2).Value = 4024
$ws.Cells.Item(57, 3).Value = 284
$ws.Cells.Item(57, 4).Value = 1539
$ws.Cells.Item(57, 5).Value = 2459
$ws.Cells.Item(57, 6).Value = 66
$ws.Cells.Item(57, 7).Value = 2
$ws.Cells.Item(57, 8).Value = 26

$ws.Cells.Item(58, 1).Value = "Argelia"
$ws.Cells.Item(58, 2).Value = 3848
$ws.Cells.Item(58, 3).Value = 0
$ws.Cells.Item(58, 4).Value = 1702
$ws.Cells.Item(58, 5).Value = 1702
$ws.Cells.Item(58, 6).Value = 22
$ws.Cells.Item(58, 7).Value = 0
$ws.Cells.Item(58, 8).Value = 444

$ws.Cells.Item(59, 1).Value = "Moldavia"
$ws.Cells.Item(59, 2).Value = 3771
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(59, 4).Value = 1182
$ws.Cells.Item(59, 5).Value = 2473
$ws.Cells.Item(59, 6).Value = 237
$ws.Cells.Item(59, 7).Value = 5
$ws.Cells.Item(59, 8).Value = 116

$ws.Cells.Item(60, 1).Value = "Luxemburgo"
$ws.Cells.Item(60, 2).Value = 3769
$ws.Cells.Item(60, 3).Value = 0
$ws.Cells.Item(60, 4).Value = 3134
$ws.Cells.Item(60, 5).Value = 546
$ws.Cells.Item(60, 6).Value = 21
$ws.Cells.Item(60, 7).Value = 0
$ws.Cells.Item(60, 8).Value = 89

$ws.Cells.Item(61, 1).Value = "Kazajistan"
$ws.Cells.Item(61, 2).Value = 3356
$ws.Cells.Item(61, 3).Value = 218
$ws.Cells.Item(61, 5).Value = 2512

# --- Libano (row 99): refresh totals ---
$ws.Cells.Item(99, 2).Value = 725
$ws.Cells.Item(99, 3).Value = 4
$ws.Cells.Item(99, 5).Value = 551

# --- Rows 144-146: Sierra Leona moves up ahead of Camboya ---
$ws.Cells.Item(144, 1).Value = "Guayana Francesa"
$ws.Cells.Item(144, 2).Value = 126
$ws.Cells.Item(144, 3).Value = 1
$ws.Cells.Item(144, 4).Value = 94
$ws.Cells.Item(144, 6).Value = 2

$ws.Cells.Item(145, 1).Value = "Sierra Leona"
$ws.Cells.Item(145, 2).Value = 124
$ws.Cells.Item(145, 3).Value = 20
$ws.Cells.Item(145, 4).Value = 21
$ws.Cells.Item(145, 5).Value = 96
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 3
$ws.Cells.Item(145, 8).Value = 7

$ws.Cells.Item(146, 1).Value = "Camboya"
$ws.Cells.Item(146, 2).Value = 122
$ws.Cells.Item(146, 3).Value = 0
$ws.Cells.Item(146, 4).Value = 119
$ws.Cells.Item(146, 5).Value = 3
$ws.Cells.Item(146, 6).Value = 1
$ws.Cells.Item(146, 8).Value = 0

# --- Guyana (row 159): refresh critical-cases count ---
$ws.Cells.Item(159, 6).Value = 3
